$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.837.39'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.887.90'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'0.7541"
$ws.Range('E5').Value = '  -2.98%  '
$ws.Range('D6').Value = "'242.40"
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('D7').Value = "'1.0000"
$ws.Range('D8').Value = "'0.3121"
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'25.36"
$ws.Range('E9').Value = '  -1.26%  '
$ws.Range('D10').Value = "'0.07121"
$ws.Range('E10').Value = '  -3.29%  '
$ws.Range('D11').Value = "'0.08482"
$ws.Range('E11').Value = '  +4.92%  '
$ws.Range('D12').Value = "'0.7603"
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('D13').Value = '1.880.09'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').Value = "'5.362"
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').Value = "'93.40"
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = "'6.134"
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '29.829.58'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = "'13.71"
$ws.Range('E18').Value = '  -1.69%  '
$ws.Range('D19').Value = "'243.52"
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('D20').Value = "'0.000007807"
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '2.141.02'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('D23').Value = "'8.002"
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('D24').Value = "'1.000"
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').Value = "'0.1600"
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').Value = "'9.380"
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').Value = "'162.80"
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').Value = "'18.73"
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').Value = "'2.031"
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').Value = "'1.486"
$ws.Range('E30').Value = '  +3.57%  '
$ws.Range('D31').Value = "'1.531"
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('D32').Value = "'4.515"
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = "'4.133"
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('D34').Value = "'0.05425"
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('D35').Value = "'1.243"
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').Value = "'0.7508"
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('D38').Value = "'2.710"
$ws.Range('E38').Value = '  +0.99%  '
$ws.Range('D39').Value = "'0.01947"
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('D40').Value = "'2.771"
$ws.Range('E40').Value = '  -1.01%  '
$ws.Range('D41').Value = "'0.4462"
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').Value = "'6.103"
$ws.Range('E42').Value = '  +2.27%  '
$ws.Range('D43').Value = '1.092.23'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('D44').Value = "'72.60"
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('D45').Value = "'0.8608"
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = "'7.722"
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('D48').Value = "'102.36"
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').Value = "'1.859"
$ws.Range('E49').Value = '  -1.51%  '
$ws.Range('D50').Value = "'3.057"
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').Value = '2.034.22'
$ws.Range('E51').Value = '  -0.61%  '
